$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.618.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.51"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4763"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2929"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06535"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.06"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7422"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.83"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.876.69"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.217"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.03"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.716.24"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007526"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.122.77"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.261"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.207"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.201"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.916"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09850"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.291"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.120"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04835"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6965"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01881"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.763"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.270"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.38"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.995"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4244"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8387"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.332"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.042"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.49"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "909.92"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3918"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.21%  "

Write-Host "Applied all updates"